$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 3")

# Row 13: clear the "x" mark in I13 and change count in J13 from 3 to 2
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = 2

# Row 14 (entry 8): fill in the previously empty time-log entry
$ws.Range("B14").Value = 43878
$ws.Range("C14").Value = 0.625
$ws.Range("D14").Value = 0.8125
$ws.Range("E14").Value = 80
$ws.Range("F14").Value = 190
$ws.Range("G14").Value = "Kodutoo MVC+konspekti korrastus"
$ws.Range("J14").Value = 1

# Row 15 (entry 9): fill in the previously empty time-log entry
$ws.Range("B15").Value = 43878
$ws.Range("C15").Value = 0.875
$ws.Range("D15").Value = 0.91666666666666663
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = 60
$ws.Range("G15").Value = "Kodutoo MVC"
$ws.Range("I15").Value = "x"
$ws.Range("J15").Value = 1

# Recalculate so the SUM(F7:F18) cached value in F19 is refreshed
$excel.Calculate()

# Update the active selection to match the saved workbook state
$ws.Range("G17").Select()
